$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 85.66330624518065
$ws.Range("C2").Value = 60.660019053970352
$ws.Range("D2").Value = 49.234431821726488
$ws.Range("E2").Value = 57.114956082693503

# Row 3 data values
$ws.Range("B3").Value = 69.666806315350954
$ws.Range("C3").Value = 42.657777568082231
$ws.Range("D3").Value = 32.907222192793327
$ws.Range("E3").Value = 57.136207392154937

# Update selection range
$ws.Range("B1:E3").Select()
